$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarterly
# data (old D:K) right to F:M.
$ws.Columns("D:E").Insert()

# The newly inserted columns come in with default/general formatting;
# copy the number formats from column F (which now holds the data that
# used to be in column D) onto the two new columns so dates/numbers
# keep rendering correctly.
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new columns with the newest quarter (D) and the
# quarter before it (E) for every data row on the sheet.
$ws.Range("D7").Value2 = 43496
$ws.Range("E7").Value2 = 43404
$ws.Range("D8").Value2 = 628100
$ws.Range("E8").Value2 = 651500
$ws.Range("D9").Value2 = 209600
$ws.Range("E9").Value2 = 221400
$ws.Range("D10").Value2 = 418500
$ws.Range("E10").Value2 = 430100
$ws.Range("D12").Value2 = 21000
$ws.Range("E12").Value2 = 22600
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = "NA"
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 36600
$ws.Range("E15").Value2 = 36200
$ws.Range("D17").Value2 = 517200
$ws.Range("E17").Value2 = 528800
$ws.Range("D18").Value2 = 110900
$ws.Range("E18").Value2 = 122700
$ws.Range("D20").Value2 = 1100
$ws.Range("E20").Value2 = 12800
$ws.Range("D21").Value2 = 180800
$ws.Range("E21").Value2 = 206000
$ws.Range("D22").Value2 = 18200
$ws.Range("E22").Value2 = 22800
$ws.Range("D23").Value2 = 93800
$ws.Range("E23").Value2 = 112700
$ws.Range("D24").Value2 = 9300
$ws.Range("E24").Value2 = -5900
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 84500
$ws.Range("E26").Value2 = 118600
$ws.Range("D27").Value2 = 84500
$ws.Range("E27").Value2 = 118600
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 18700
$ws.Range("E29").Value2 = -17900
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -1100
$ws.Range("E32").Value2 = -12800
$ws.Range("D33").Value2 = 103200
$ws.Range("E33").Value2 = 100700
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 103200
$ws.Range("E35").Value2 = 100700
$ws.Range("D38").Value2 = 43496
$ws.Range("E38").Value2 = 43404
$ws.Range("D41").Value2 = 146600
$ws.Range("E41").Value2 = 77700
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 389700
$ws.Range("E43").Value2 = 374700
$ws.Range("D44").Value2 = 486200
$ws.Range("E44").Value2 = 468800
$ws.Range("D45").Value2 = 170500
$ws.Range("E45").Value2 = 169700
$ws.Range("D46").Value2 = 1193000
$ws.Range("E46").Value2 = 1090900
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 1023800
$ws.Range("E48").Value2 = 976000
$ws.Range("D49").Value2 = 3963400
$ws.Range("E49").Value2 = 3913400
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 118400
$ws.Range("E52").Value2 = 132500
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 6298600
$ws.Range("E54").Value2 = 6112800
$ws.Range("D57").Value2 = 129600
$ws.Range("E57").Value2 = 146400
$ws.Range("D58").Value2 = 445500
$ws.Range("E58").Value2 = 37100
$ws.Range("D59").Value2 = 363500
$ws.Range("E59").Value2 = 353000
$ws.Range("D60").Value2 = 938600
$ws.Range("E60").Value2 = 536500
$ws.Range("D61").Value2 = 1686900
$ws.Range("E61").Value2 = 1985700
$ws.Range("D62").Value2 = 247600
$ws.Range("E62").Value2 = 282800
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 2873300
$ws.Range("E66").Value2 = 2805200
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 2664400
$ws.Range("E72").Value2 = 2576000
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 3425300
$ws.Range("E76").Value2 = 3307600
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43496
$ws.Range("E80").Value2 = 43404
$ws.Range("D81").Value2 = 103200
$ws.Range("E81").Value2 = 100700
$ws.Range("D83").Value2 = 68800
$ws.Range("E83").Value2 = 70500
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 101800
$ws.Range("E89").Value2 = 236600
$ws.Range("D91").Value2 = -79200
$ws.Range("E91").Value2 = -43400
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -129200
$ws.Range("E94").Value2 = -46500
$ws.Range("D96").Value2 = 0
$ws.Range("E96").Value2 = -1400
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = 93200
$ws.Range("E100").Value2 = -263300
$ws.Range("D101").Value2 = 900
$ws.Range("E101").Value2 = -2200
$ws.Range("D102").Value2 = 66700
$ws.Range("E102").Value2 = -75400
